$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values that changed ---
$ws.Range("C2").Value = -46
$ws.Range("D2").Value = 0.999
$ws.Range("G2").Value = "129.733 mA"
$ws.Range("H2").Value = "20.135 mA"
$ws.Range("I2").Value = "0.007 mA"
$ws.Range("M2").Value = 43398.388136574074

# --- Create row 3 by copying formatting from row 2, then set its values ---
$ws.Range("A2:M2").Copy($ws.Range("A3:M3"))

$ws.Range("A3").Value = "PASS"
$ws.Range("B3").Value = "000B57FFFEF609E6"
$ws.Range("C3").Value = -55
$ws.Range("D3").Value = 0.99199999999999999
$ws.Range("E3").Value = -55
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "121.584 mA"
$ws.Range("H3").Value = "20.332 mA"
$ws.Range("I3").Value = "0.013 mA"
$ws.Range("J3").Value = "OK"
$ws.Range("K3").Value = "OK"
$ws.Range("L3").Value = "0x0000"
$ws.Range("M3").Value = 43398.572187500002

# --- Update the active selection to H14 ---
$ws.Range("H14").Select()
